# Översikt HAMMARÖ.xlsx update
#
# The commit reorders the data rows (rows 2-33) of the sheet and bumps the
# "Förändrad" (column C) timestamp on every row from 2026-02-26 (46079) to
# 2026-02-28 (46081). No data is actually added or removed - every row's
# full content (A..Z) simply moves to a new row position; a handful of rows
# stay put. We reproduce this by:
#   1. snapshotting every existing data row (keyed by its "Beteckning" / A
#      column id, which is unique per row) including formulas vs values,
#   2. clearing the data rows,
#   3. writing each row back out in the new order, refreshing column C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 33
$lastCol = 26   # column Z
$newForandradSerial = 46081   # 2026-02-28, was 46079 (2026-02-26)

# Final row order, by the value of column A ("Beteckning"). Derived from the
# target XML: most rows keep their neighbours but several pairs/triples swap
# places (e.g. old row 3 <-> old row 4, old rows 9-12 rotate by one, etc).
$newOrder = @(
    "A 68622-2021",
    "A 47571-2025",
    "A 61380-2023",
    "A 3573-2026",
    "A 24616-2022",
    "A 45407-2025",
    "A 53276-2023",
    "A 58382-2024",
    "A 55068-2023",
    "A 20755-2021",
    "A 24618-2022",
    "A 16890-2021",
    "A 56835-2021",
    "A 51008-2023",
    "A 57893-2023",
    "A 37934-2024",
    "A 56799-2022",
    "A 45423-2025",
    "A 41661-2025",
    "A 41895-2023",
    "A 41899-2023",
    "A 55069-2023",
    "A 58383-2024",
    "A 53369-2024",
    "A 13356-2022",
    "A 42951-2023",
    "A 42955-2023",
    "A 42960-2023",
    "A 13354-2022",
    "A 55066-2023",
    "A 2769-2023",
    "A 42957-2023"
)

# --- 1. snapshot every row, keyed by its id -------------------------------
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $id = $ws.Cells.Item($r, 1).Value2()
    $rowData = @{}
    for ($c = 1; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($cell.HasFormula()) {
            $rowData[$c] = @{ kind = "formula"; val = $cell.Formula() }
        } else {
            $v = $cell.Value2()
            if ($v -eq $null) {
                $rowData[$c] = @{ kind = "empty"; val = $null }
            } else {
                $rowData[$c] = @{ kind = "value"; val = $v }
            }
        }
    }
    $snapshot[$id] = $rowData
}

# --- 2. clear the data rows (content only, formatting is shared per-column
#         and stays put) ----------------------------------------------------
$clearRange = $ws.Range($ws.Cells.Item($firstRow, 1), $ws.Cells.Item($lastRow, $lastCol))
$clearRange.ClearContents()

# --- 3. write every row back out in its new position ----------------------
$destRow = $firstRow
foreach ($id in $newOrder) {
    $rowData = $snapshot[$id]
    for ($c = 1; $c -le $lastCol; $c++) {
        $entry = $rowData[$c]
        if ($entry.kind -eq "formula") {
            $ws.Cells.Item($destRow, $c).Formula = $entry.val
        } elseif ($entry.kind -eq "value") {
            $ws.Cells.Item($destRow, $c).Value = $entry.val
        }
        # "empty" entries: nothing to write, row was already cleared.
    }
    # column C ("Förändrad") is bumped to the new date on every row.
    $ws.Cells.Item($destRow, 3).Value = $newForandradSerial
    $destRow++
}
